$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.2881169905109251;  C = 0.04103571897497393; D = 3.223369029078222;  E = 0.5333859586016987;  G = 4.085907697165819 }
    3 = @{ B = 0.04172184405617529; C = 0.3048912486333797;  D = 0.7210945179870265; E = 0.5333859586016987;  G = 1.60109356927828 }
    4 = @{ B = 0.1169995834814548;  C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 2.998467759612273 }
    5 = @{ B = 3.272327238179451;   C = 1.626987699542094;   D = 0.1496068669990043; E = 0.5333859586016987;  G = 5.582307763322248 }
    6 = @{ B = 3.272327238179451;   C = 9.983522426115931;   D = 189.6080260415259;  E = 13.86384647080068;   G = 216.727722176622 }
    7 = @{ B = 3.272327238179451;   C = 1.626987699542094;   D = 3.223369029078222;  E = 0.5333859586016987;  G = 8.656069925401464 }
    8 = @{ B = 0.6545652718822623; C = 1.626987699542094;   D = 0.7210945179870265; E = 0.5333859586016987;  G = 3.536033448013082 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
